# CMP73010-Ass1.docx edit:
#   - Paragraph 2 ("This is the document you will need to change...") is
#     replaced with the merge-request instructions, split across four runs
#     (matching the original author's multi-run typing/paste pattern).
#   - Paragraph 3 ("Add instructions for your tutor...") has all of its
#     text removed but keeps the existing "_GoBack" bookmark pair.
#   - A new, final empty paragraph is added before the section break.

$d = $word.ActiveDocument

# --- Paragraph 2: replace the placeholder sentence with the merge how-to ---
$para2Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>To merge a file with</w:t></w:r><w:r><w:t xml:space="preserve"> mainline, you need to login to your git hub account and you need to pull request a document after that you can click “merge</w:t></w:r><w:r><w:t xml:space="preserve"> pull request</w:t></w:r><w:r><w:t xml:space="preserve">” button at the bottom of your pull request to merge your change and click on “Confirm Merge” to merge the changes. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $d.Paragraphs(2).Range.InsertXML($para2Xml)

# --- Paragraph 3: drop the text but keep the "_GoBack" bookmark pair ---
# (InsertXML on the whole paragraph range with no runs leaves the bookmark
#  in place and mints a fresh trailing empty paragraph, which is exactly
#  the extra blank paragraph the edit adds before the section break.)
$para3Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$null = $d.Paragraphs(3).Range.InsertXML($para3Xml)

Write-Output "Edit applied. Paragraph count: $($d.Paragraphs.Count)"
